# Update the "orders" sheet: rename the "Order" / "Order Number" columns
# to "Order_Nr" / "Order_Number" (s3 file name update), and refresh the
# active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("orders")

# Write C2 first so "Order_Number" lands before "Order_Nr" in the shared
# string table (matches upstream ordering). The leading apostrophe forces
# Excel's "stored as text" quote-prefix formatting on this cell, same as
# the authored workbook.
$ws.Range("C2").Value = "'Order_Number"
$ws.Range("A2").Value = "Order_Nr"

# Move the selection/active cell on sheet "orders" from B7 to A2.
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
